$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7901163333333333
$ws.Range("H2").Value = 2.370349
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 224.2321046666667
$ws.Range("N2").Value = 672.696314
$ws.Range("O2").Value = 0.9009864013525987
$ws.Range("P2").Value = 0.9009864013525988
$ws.Range("Q2").Value = 177.1694483548429
$ws.Range("R2").Value = 1594.525035193586
$ws.Range("S2").Value = 0.9009864013525987
$ws.Range("T2").Value = 0.9009864013525988

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7901163333333333
$ws.Range("H3").Value = 2.370349
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.023704333333333
$ws.Range("N3").Value = 3.071113
$ws.Range("O3").Value = 0.004113343558497904
$ws.Range("P3").Value = 0.004113343558497904
$ws.Range("Q3").Value = 0.8088455142707778
$ws.Range("R3").Value = 7.279609628437
$ws.Range("S3").Value = 0.004113343558497904
$ws.Range("T3").Value = 0.004113343558497904

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7901163333333333
$ws.Range("H4").Value = 2.370349
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.61820766666667
$ws.Range("N4").Value = 70.854623
$ws.Range("O4").Value = 0.0949002550889034
$ws.Range("P4").Value = 0.09490025508890343
$ws.Range("Q4").Value = 18.66113164149189
$ws.Range("R4").Value = 167.950184773427
$ws.Range("S4").Value = 0.0949002550889034
$ws.Range("T4").Value = 0.09490025508890343
